$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Purchase 22-23" (sheet1)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Purchase 22-23")

# Drop the "Collective Trade Links Pvt Ltd" entry completely: its two
# invoice lines (rows 5-6) plus the blank separator row below them (row 7).
$ws1.Range("A5:A7").EntireRow.Delete()

# Drop the first "68/23-24" invoice line of the Namrata Rubber entry
# (old row 2), leaving the "71/23-24" line (old row 3) on its own.
$ws1.Range("A2").EntireRow.Delete()

# The surviving Namrata Rubber line is now row 2 on its own: restore its
# Sr. No and make the outstanding total just "=E2" (no second addend).
$ws1.Range("A2").Value = 1
$ws1.Range("F2").Formula = "=E2"

# Renumber the Sr. No column for the remaining entries sequentially.
$ws1.Range("A4").Value = 2
$ws1.Range("A6").Value = 3
$ws1.Range("A8").Value = 4
$ws1.Range("A10").Value = 5
$ws1.Range("A12").Value = 6
$ws1.Range("A14").Value = 7

# Append the new "Asha Enterprises" entry as rows 16-17 (entry #8); row 15
# stays as the usual blank separator. Bring over matching cell formatting
# from existing rows first, then fill in the values/formula.
$ws1.Range("A8:E8").Copy() | Out-Null
$ws1.Range("A16:E16").PasteSpecial(-4122) | Out-Null
$ws1.Range("A17:E17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("C8").Copy() | Out-Null
$ws1.Range("F16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("F8").Copy() | Out-Null
$ws1.Range("F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("A16").Value = 8
$ws1.Range("B16").Value = 45254
$ws1.Range("C16").Value = 2067
$ws1.Range("D16").Value = "Asha Enterprises"
$ws1.Range("E16").Value = 2561

$ws1.Range("B17").Value = 45262
$ws1.Range("C17").Value = 2074
$ws1.Range("D17").Value = "Asha Enterprises"
$ws1.Range("E17").Value = 496
$ws1.Range("F17").Formula = "=E16+E17"

# Sheet view: selection moves, and the sheet is no longer the active tab.
$ws1.Range("D31").Select()

# ------------------------------------------------------------------
# Sheet "Sale 22-23" (sheet2)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# Drop the stray one-off formula row at the bottom of the sheet.
$ws2.Range("A27").EntireRow.Delete()

# This sheet becomes the active tab/sheet, with a new selection.
$ws2.Activate()
$ws2.Range("K7").Select()
